$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")

# Updated metrics for existing classifiers (kNN, SVM, LR, NB) - rows 2-5
$row2 = @(0.5154639175257731,0.5501022163793501,0.5262617327632807,0.4807650635718315,0.5670103092783505,0.632590446989868,0.5729038281979457,0.5441641680825956,0.5670103092783505,0.6212897563171373,0.5693567251461988,0.5577232950839509,0.6082474226804123,0.6126692698575859,0.6113813946631284,0.5999856053719078,0.468041237113402,0.5370471273372495,0.468391567153177,0.4608091376274689,0.468041237113402,0.5370471273372495,0.468391567153177,0.4608091376274689)
$row3 = @(0.6948453608247424,0.7238971077299784,0.6940940586761021,0.6912806583050726,0.7092783505154638,0.7373028001315965,0.707936507936508,0.7084164375674769,0.7092783505154638,0.7373028001315965,0.707936507936508,0.7084164375674769,0.6927835051546392,0.7163858237984029,0.6918688879060396,0.6877167398491423,0.7319587628865979,0.7563400033527428,0.7302452208953758,0.726973369957835,0.7319587628865979,0.7563400033527428,0.7302452208953758,0.726973369957835)
$row4 = @(0.6907216494845361,0.7115741915905511,0.6895660720428524,0.6864554750414388,0.7278350515463917,0.7531931742185816,0.7258115877930119,0.7283094043839362,0.7278350515463917,0.7531931742185816,0.7258115877930119,0.7283094043839362,0.7010309278350516,0.7242188254153156,0.697544842498403,0.6938546586702822,0.7298969072164948,0.7576090174642015,0.726889773453241,0.7267085276261819,0.7298969072164948,0.7576090174642015,0.726889773453241,0.7267085276261819)
$row5 = @(0.6556701030927835,0.6688792251107503,0.652435991940636,0.6441128945624456,0.7278350515463917,0.7435934855381079,0.7235824856258292,0.7217225062629988,0,0,0,0,0.6536082474226804,0.6855695759289765,0.6482008943928449,0.6359419755845588,0.6989690721649485,0.7414403171431545,0.6954818418595508,0.6920731469137873,0,0,0,0)

# New row 6 - Ensemble classifier (first approach)
$row6 = @(0.6927835051546392,0.7148699056859644,0.6954189395056268,0.6884880347501156,0.7278350515463917,0.7645850209996984,0.7285566858322277,0.7282455651959145,0.7154639175257731,0.7468185316526796,0.7154906875030714,0.7133057519877879,0.6969072164948453,0.7118356590354301,0.6967310432945107,0.6927819415178955,0.7443298969072166,0.7666078794697349,0.7428876111848248,0.742514809634382,0.7154639175257731,0.7468185316526796,0.7154906875030714,0.7133057519877879)

$data = @{ 2 = $row2; 3 = $row3; 4 = $row4; 5 = $row5; 6 = $row6 }

foreach ($r in @(2,3,4,5,6)) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}

# Row 6 label cell (A6) - copy formatting from A5 (bold, centered, bordered) then set its text
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Ensemble"
